$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in rows 8, 10, 11, 12, 13, 14, 15
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("C9").Value = 16

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10

$ws.Range("D13").Value = 8

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# Add two new rows (16 and 17)
$ws.Range("A16").Value = 14
$ws.Range("A16").Style = $ws.Range("A15").Style
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("A17").Style = $ws.Range("A15").Style
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
